# Payment report template update + dropdown authorities edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the column header text in G6 (service-invoice number -> debt-invoice number)
$ws.Range("G6").Value = "เลขที่ใบแจ้งค่าหนี้"

# Reset the scrolled view back to the top-left and move the active
# selection to H10, matching the saved view state in the workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H10").Select()

# Shrink the saved workbook window height.
$excel.ActiveWindow.Height = 5385
